$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.318907260894775
$ws.Range("B1").Value = 2.263653516769409
$ws.Range("C1").Value = 1.616185903549194
$ws.Range("D1").Value = 1.469524264335632
$ws.Range("E1").Value = 1.457105040550232
